$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.33"
$ws.Range("E2").Value = "'-3.19%"
$ws.Range("D3").Value = "'40.98"
$ws.Range("E3").Value = "'-2.22%"
$ws.Range("E4").Value = "'-3.00%"
$ws.Range("D5").Value = "'0.07622"
$ws.Range("E5").Value = "'-5.59%"
$ws.Range("D6").Value = "'4.234"
$ws.Range("E6").Value = "'-3.28%"
$ws.Range("E7").Value = "'-8.96%"
$ws.Range("D8").Value = "'0.9066"
$ws.Range("E8").Value = "'-2.31%"
$ws.Range("D9").Value = "'0.09886"
$ws.Range("E9").Value = "'-12.04%"
$ws.Range("D10").Value = "'0.1767"
$ws.Range("E10").Value = "'-4.78%"
$ws.Range("D11").Value = "'0.09126"
$ws.Range("E11").Value = "'-1.41%"
$ws.Range("D12").Value = "'0.04408"
$ws.Range("E12").Value = "'-3.39%"
$ws.Range("E13").Value = "'-0.09%"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("E14").Value = "'-2.52%"
$ws.Range("D15").Value = "'0.005877"
$ws.Range("E15").Value = "'-1.60%"
$ws.Range("D16").Value = "'3.369"
$ws.Range("E16").Value = "'0.41%"
$ws.Range("D17").Value = "'2.419"
$ws.Range("E17").Value = "'-5.03%"
$ws.Range("E18").Value = "'-2.95%"
$ws.Range("D19").Value = "'6.832"
$ws.Range("E19").Value = "'-7.31%"
$ws.Range("D20").Value = "'0.1348"
$ws.Range("E20").Value = "'-2.51%"
$ws.Range("D21").Value = "'0.2841"
$ws.Range("E21").Value = "'9.06%"
$ws.Range("D22").Value = "'0.04165"
$ws.Range("E22").Value = "'-0.20%"
$ws.Range("D23").Value = "'0.001214"
$ws.Range("E23").Value = "'-2.39%"
$ws.Range("D24").Value = "'0.004068"
$ws.Range("E24").Value = "'-5.75%"
$ws.Range("D25").Value = "'0.0001300"
$ws.Range("E25").Value = "'6.40%"
$ws.Range("D26").Value = "'0.0003005"
$ws.Range("E26").Value = "'0.60%"
$ws.Range("D38").Value = "'0.02418"
$ws.Range("E38").Value = "'-6.05%"
$ws.Range("D39").Value = "'0.05140"
$ws.Range("E39").Value = "'-5.52%"
$ws.Range("E40").Value = "'-3.03%"
$ws.Range("D41").Value = "'0.1307"
$ws.Range("E41").Value = "'-6.08%"
$ws.Range("D42").Value = "'0.007060"
$ws.Range("E42").Value = "'-6.57%"
$ws.Range("D43").Value = "'0.001948"
$ws.Range("E43").Value = "'-6.50%"
$ws.Range("D44").Value = "'0.008372"
$ws.Range("E44").Value = "'1.59%"
$ws.Range("D45").Value = "'0.3049"
$ws.Range("E45").Value = "'-3.12%"
$ws.Range("D46").Value = "'0.00006376"
$ws.Range("E46").Value = "'-6.25%"
$ws.Range("E47").Value = "'-0.24%"
$ws.Range("E48").Value = "'-26.95%"
$ws.Range("D49").Value = "'0.007430"
$ws.Range("E49").Value = "'119.12%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.24%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.24%"
